# Append two new daily rows (date + value) to each of the 7 worksheets,
# matching the layout of the existing data (date in column A formatted as
# a date/time, value in column B). The final new row is a zero placeholder
# for the latest date.

$wb = $excel.ActiveWorkbook

$newDate1 = 45967   # 2025-11-06
$newDate2 = 45968   # 2025-11-07
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# New value for the 45967 row, per worksheet (1-indexed by tab order).
$newValues = @(3128, 1171, 1256, 1860, 753, 1497, 2963)

for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $row1 = 104
    $row2 = 105

    $ws.Cells.Item($row1, 1).Value = $newDate1
    $ws.Cells.Item($row1, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row1, 2).Value = $newValues[$i - 1]

    $ws.Cells.Item($row2, 1).Value = $newDate2
    $ws.Cells.Item($row2, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row2, 2).Value = 0
}
